$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range('D2').NumberFormat = '@'
$ws.Range('D2').Value = '247.24'
$ws.Range('D3').NumberFormat = '@'
$ws.Range('D3').Value = '21.76'
$ws.Range('D4').NumberFormat = '@'
$ws.Range('D4').Value = '5.444'
$ws.Range('D5').NumberFormat = '@'
$ws.Range('D5').Value = '0.05669'
$ws.Range('D6').NumberFormat = '@'
$ws.Range('D6').Value = '3.371'
$ws.Range('D7').NumberFormat = '@'
$ws.Range('D7').Value = '0.8014'
$ws.Range('D8').NumberFormat = '@'
$ws.Range('D8').Value = '1.035'
$ws.Range('B9').Value = 'WazirX'
$ws.Range('C9').Value = 'https://coinranking.com/coin/6QK-8hUZ+wazirx-wrx'
$ws.Range('D9').NumberFormat = '@'
$ws.Range('D9').Value = '0.1431'
$ws.Range('E9').Value = '8WazirXWRX'
$ws.Range('B10').Value = 'MandalaExchangeToken'
$ws.Range('C10').Value = 'https://coinranking.com/coin/lviNIbma2Xuqs+mandalaexchangetoken-mdx'
$ws.Range('D10').NumberFormat = '@'
$ws.Range('D10').Value = '0.07238'
$ws.Range('E10').Value = '9MandalaExchangeTokenMDX'
$ws.Range('B11').Value = 'LiechtensteinCryptoassetsExchange'
$ws.Range('C11').Value = 'https://coinranking.com/coin/v4IW9oaF+liechtensteincryptoassetsexchange-lcx'
$ws.Range('D11').NumberFormat = '@'
$ws.Range('D11').Value = '0.03157'
$ws.Range('E11').Value = '10LiechtensteinCryptoassetsExchangeLCX'
$ws.Range('B12').Value = 'BitrueCoin'
$ws.Range('C12').Value = 'https://coinranking.com/coin/SLYjzF4ty+bitruecoin-btr'
$ws.Range('D12').NumberFormat = '@'
$ws.Range('D12').Value = '0.02944'
$ws.Range('E12').Value = '11BitrueCoinBTR'
$ws.Range('B13').Value = 'BitMartToken'
$ws.Range('C13').Value = 'https://coinranking.com/coin/6uzcPMFgWUJNH+bitmarttoken-bmx'
$ws.Range('D13').NumberFormat = '@'
$ws.Range('D13').Value = '0.09277'
$ws.Range('E13').Value = '12BitMartTokenBMX'
$ws.Range('B14').Value = 'BitForexToken'
$ws.Range('C14').Value = 'https://coinranking.com/coin/2nh5ugplNocUp+bitforextoken-bf'
$ws.Range('D14').NumberFormat = '@'
$ws.Range('D14').Value = '0.001652'
$ws.Range('E14').Value = '13BitForexTokenBF'
$ws.Range('B15').Value = 'MCDex'
$ws.Range('C15').Value = 'https://coinranking.com/coin/3nMM61qeg+mcdex-mcb'
$ws.Range('D15').NumberFormat = '@'
$ws.Range('D15').Value = '3.220'
$ws.Range('E15').Value = '14MCDexMCB'
$ws.Range('B16').Value = 'CoinExToken'
$ws.Range('C16').Value = 'https://coinranking.com/coin/APDVU0XEViZ2o+coinextoken-cet'
$ws.Range('D16').NumberFormat = '@'
$ws.Range('D16').Value = '0.04726'
$ws.Range('E16').Value = '15CoinExTokenCET'
$ws.Range('B17').Value = 'One'
$ws.Range('C17').Value = 'https://coinranking.com/coin/6Lga5NiXX3rT+one-one'
$ws.Range('D17').NumberFormat = '@'
$ws.Range('D17').Value = '0.0005900'
$ws.Range('E17').Value = '16OneONE'
$ws.Range('D18').NumberFormat = '@'
$ws.Range('D18').Value = '0.006480'
$ws.Range('D19').NumberFormat = '@'
$ws.Range('D19').Value = '0.005021'
$ws.Range('E19').Value = '18HotbitTokenHTBBestin24h'
$ws.Range('D21').NumberFormat = '@'
$ws.Range('D21').Value = '0.0001500'
$ws.Range('D22').NumberFormat = '@'
$ws.Range('D22').Value = '0.0003200'
$ws.Range('D24').NumberFormat = '@'
$ws.Range('D24').Value = '6.429'
$ws.Range('D25').NumberFormat = '@'
$ws.Range('D25').Value = '2.126'
$ws.Range('D40').NumberFormat = '@'
$ws.Range('D40').Value = '0.04082'
$ws.Range('B41').Value = 'BKEXToken'
$ws.Range('C41').Value = 'https://coinranking.com/coin/IPeThtYgk+bkextoken-bkk'
$ws.Range('D41').NumberFormat = '@'
$ws.Range('D41').Value = '0.1041'
$ws.Range('E41').Value = '40BKEXTokenBKK'
$ws.Range('B42').Value = 'CEJI'
$ws.Range('C42').Value = 'https://coinranking.com/coin/SbKjCVJCh+ceji-ceji'
$ws.Range('D42').NumberFormat = '@'
$ws.Range('D42').Value = '0.003201'
$ws.Range('E42').Value = '41CEJICEJI'
$ws.Range('B43').Value = 'KickToken'
$ws.Range('C43').Value = 'https://coinranking.com/coin/F_Yv9Cu7pPL3Y+kicktoken-kick'
$ws.Range('D43').NumberFormat = '@'
$ws.Range('D43').Value = '0.003255'
$ws.Range('E43').Value = '42KickTokenKICK'
$ws.Range('D44').NumberFormat = '@'
$ws.Range('D44').Value = '0.009081'
$ws.Range('D45').NumberFormat = '@'
$ws.Range('D45').Value = '0.00005847'
$ws.Range('D46').NumberFormat = '@'
$ws.Range('D46').Value = '0.00000000750'
$ws.Range('D47').NumberFormat = '@'
$ws.Range('D47').Value = '0.7852'
$ws.Range('D48').NumberFormat = '@'
$ws.Range('D48').Value = '0.01185'
$ws.Range('D49').NumberFormat = '@'
$ws.Range('D49').Value = '0.00002100'
$ws.Range('D50').NumberFormat = '@'
$ws.Range('D50').Value = '0.01010'
